$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename DataID -> DataId, add Duration (F) and CoolTime (G) headers
$ws.Range("A1").Value = "DataId"
$ws.Range("F1").Value = "Duration"
$ws.Range("G1").Value = "CoolTime"

# Row 2 (Skill_Boss_Dash / 350): Value changed 10 -> 200, add Duration=2, CoolTime=15
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 15

# Row 3 (Speedchange_Down / 400): add Duration=5, CoolTime=0
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 0

# Row 4 (Jumpchange_Up / 401): add Duration=0, CoolTime=0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Move the active selection to match the saved view state
$ws.Range("H13").Select()
